$d = $word.ActiveDocument

# 1. Fix the typo "ti returns" -> "it returns" in the get_mailbox_for_pid paragraph.
$d.Content.Find.Execute("when creating the mailbox, ti returns an error.", $true, $false, $false, $false, $false, $true, 1, $false, "when creating the mailbox, it returns an error.", 2) | Out-Null

# 2. Append the new "About sigint" section at the end of the document.

# -- three blank paragraphs
$p = $d.Paragraphs.Last
$r = $p.Range
$r.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$r = $p.Range
$r.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$r = $p.Range
$r.InsertParagraphAfter()

# -- heading paragraph "About sigint" (bold applied after the remaining paragraphs
#    are created, so the bold formatting doesn't leak into the following paragraphs)
$p = $d.Paragraphs.Last
$r = $p.Range
$r.InsertParagraphAfter()
$headingPara = $d.Paragraphs.Last
$headingPara.Range.Text = "About sigint"

# -- SIGINT explanation paragraph
$p = $d.Paragraphs.Last
$r = $p.Range
$r.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.Text = "So SIGINT seems to not call sys_exit, which is awesome because that means that our mailbox cleanup is never called. A simple way around this is to spawn a " + [char]0x201C + "garage collection" + [char]0x201D + " thread that can wake up every once and a while and check a linked list of pid_t so see if mailboxes are still alive. When we create the mailbox add a copy of the pid_t to the linked list. This could be wasteful though if there are a lot of active mailboxes. "

# -- final paragraph about signal handler / task_struct (keeps the _GoBack bookmark)
$p = $d.Paragraphs.Last
$r = $p.Range
$r.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.Text = "You could also install a signal handler to all of the processes you make a mailbox for that destroys the mailbox before handing it off to the default handler.  There" + [char]0x2019 + "s a member of task_struct for that."

# now make the heading bold, after all subsequent paragraphs already exist
$headingPara.Range.Font.Bold = 1
